$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency data (prices, volume %, and row content swaps)
# Force text number format before writing numeric-looking strings so Excel
# stores them as text (matching original inlineStr cells), then reset style
# back to Normal so no stray style index is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.788.10'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.77%  '
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.969.07'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +4.21%  '
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '483.69'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +8.81%  '
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.70'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +3.91%  '
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.627'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.52%  '
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.998'
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.734'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.11%  '
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.170'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +9.67%  '
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000358'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +11.67%  '
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.51'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.04%  '
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.599.68'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +4.73%  '
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.51'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.88%  '
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.81'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.56%  '
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.979.12'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +4.21%  '
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.13%  '
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.06'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.54%  '
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.30%  '
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.737.10'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.64%  '
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '441.22'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.95%  '
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.35'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.43%  '
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.46'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.17%  '
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.38'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.86%  '
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.63'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +6.63%  '
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '38.62'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +3.49%  '
$ws.Range("E26").Style = "Normal"

$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("B27").Style = "Normal"

$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("C27").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.03'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +6.28%  '
$ws.Range("E27").Style = "Normal"

$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = 'Filecoin'
$ws.Range("B28").Style = "Normal"

$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("C28").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.11'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +4.34%  '
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '729.13'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.26%  '
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '13.25'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.78%  '
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.128'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.75%  '
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.56%  '
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '43.01'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.60%  '
$ws.Range("E33").Style = "Normal"

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'OKB'
$ws.Range("B34").Style = "Normal"

$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("C34").Style = "Normal"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '60.23'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +4.04%  '
$ws.Range("E34").Style = "Normal"

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'PEPE'
$ws.Range("B35").Style = "Normal"

$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("C35").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0859'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +27.90%  '
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.89%  '
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.41'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.10%  '
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.12%  '
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.05'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +6.04%  '
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.92'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +9.98%  '
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.64'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +5.87%  '
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.142'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.17%  '
$ws.Range("E43").Style = "Normal"

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'TheGraph'
$ws.Range("B44").Style = "Normal"

$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("C44").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.337'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.55%  '
$ws.Range("E44").Style = "Normal"

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("B45").Style = "Normal"

$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("C45").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.12%  '
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("B46").Style = "Normal"

$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("C46").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.21'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +4.72%  '
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.41'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.12%  '
$ws.Range("E47").Style = "Normal"

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'Monero'
$ws.Range("B48").Style = "Normal"

$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("C48").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '149.02'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.34%  '
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("B49").Style = "Normal"

$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("C49").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.21'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.67%  '
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.91'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.41%  '
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '25.19'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.18%  '
$ws.Range("E51").Style = "Normal"
